# Port_Letter.xlsx template fix: insert 4 hidden spacer rows before the
# existing 0.4pt spacer block (old rows 27-30) so the "letter body" table
# has extra blank rows to grow into, and bump that block's row heights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 4 new rows before row 27 (formats copy down from row 26, and
#    every formula/range below - merges, data validations, defined names -
#    shifts down by 4 automatically).
$ws.Range('A27:A30').EntireRow.Insert()

# 2) Re-apply the row-height plan across rows 25-41 (new row numbering).
$ws.Range('A25:A26').EntireRow.RowHeight = 1.25
$ws.Range('A27:A28').EntireRow.RowHeight = 1.25
$ws.Range('A30:A41').EntireRow.RowHeight = 11.75

# 3) New spacer-row formatting: wrap text on the representative-name cell.
$ws.Range('A42').WrapText = $true

# 4) Fix up the two conditional-formatting rules, whose sqref (and, for the
#    second, the relative row in its formula) does not auto-shift with the
#    row insert.
$rule1 = $ws.Range('A24:F26').FormatConditions.Item(1)
$rule1.ModifyAppliesToRange($ws.Range('A24:F30'))

$rule2 = $ws.Range('H29:N44').FormatConditions.Item(1)
$rule2.ModifyAppliesToRange($ws.Range('H33:N48'))
$rule2.Formula1 = '=$H33<>0'

# 4b) The second date-validation's Formula1 is a relative same-row ref
#     ("M29"/"M33") that also doesn't auto-shift with the row insert.
$ws.Range('N33').Validation.Formula1 = 'M33'

# 5) Fix up the print area and the named ranges that point below row 26 -
#    these also don't auto-shift.
$ws.PageSetup.PrintArea = '$A$1:$F$58'

$shiftedNames = @{
  'Merge_end'                = '$A$47'
  'Pg_end'                   = '$F$50'
  'Seal_seller_end'          = '$E$50'
  'Seal_seller_start'        = '$D$46'
  'Sign_seller_end'          = '$E$49'
  'Sign_seller_start'        = '$D$49'
  'Выгрузка_ответственный'   = '$A$35'
  'Грузовые_борт_склад'      = '$A$37'
  'Грузовые_склад_авто'      = '$A$38'
  'Имя_представитель'        = '$A$42'
  'Исполнитель_информация'   = '$A$40'
  'Контрольный_звонок'       = '$A$45'
  'Образцы_выдача'           = '$A$31'
  'Образцы_подвал'           = '$A$32'
  'Письмо_описание_подвал'   = '$A$33'
  'Подписант'                = '$F$49'
  'Подписант_комментарий'    = '$A$49'
  'Покупатель_телефон'       = '$A$36'
  'Расходы_компания'         = '$A$34'
  'Телефон_представитель'    = '$A$43'
  'Хранение'                 = '$A$39'
}
foreach ($nm in $shiftedNames.Keys) {
    $wb.Names.Item($nm).RefersTo = '=Port_Letter!' + $shiftedNames[$nm]
}

# 6) Match the author's final selection in the sheet view.
$ws.Range('C43').Select()

Write-Output "Port_Letter: inserted 4 rows before row 27, fixed heights/styles/ranges"
